# Applies the "research paper added (2)" edit:
#  1. Re-style the "Abstract" paragraph: remove the underline and make it
#     bold + centered (matching the new section headings).
#  2. Append four new centered/bold section headings ("Introduction",
#     "Technical / scientific core sections", "Analysis, discussion,
#     conclusions", "Bibliography"), each followed by a blank paragraph.

$d = $word.ActiveDocument

$flatOpcHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$flatOpcFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function ConvertTo-FlatOpcRun($innerBodyXml) {
    return $flatOpcHeader + $innerBodyXml + $flatOpcFooter
}

function Get-HeadingParagraphXml($text) {
    return "<w:p><w:pPr><w:jc w:val=`"center`"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>$text</w:t></w:r></w:p>"
}

# --- 1. Re-style the existing "Abstract" paragraph --------------------------
$abstractParagraph = $d.Paragraphs(2)
$abstractParagraph.Range.InsertXML((ConvertTo-FlatOpcRun (Get-HeadingParagraphXml "Abstract")))

# --- 2. Append the new section headings -------------------------------------
# InsertXML is only reliable for a couple of new paragraphs at a time on this
# host, so insertions are performed in small batches: grow the document by one
# placeholder paragraph, then replace that placeholder (via InsertXML) with
# the batch's real paragraphs.
function Add-ParagraphsAtEnd($innerBodyXml) {
    $lastParagraph = $d.Paragraphs($d.Paragraphs.Count)
    $lastParagraph.Range.InsertParagraphAfter()
    $placeholder = $d.Paragraphs($d.Paragraphs.Count)
    $placeholder.Range.InsertXML((ConvertTo-FlatOpcRun $innerBodyXml))
}

Add-ParagraphsAtEnd ((Get-HeadingParagraphXml "Introduction") + "<w:p/>")
Add-ParagraphsAtEnd ((Get-HeadingParagraphXml "Technical / scientific core sections") + "<w:p/>")
Add-ParagraphsAtEnd ((Get-HeadingParagraphXml "Analysis, discussion, conclusions") + "<w:p/>")
Add-ParagraphsAtEnd (Get-HeadingParagraphXml "Bibliography")
